# Updated cryptos list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.676.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.32%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.674.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.66%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +7.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "655.51"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.421"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.84%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.05%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.673.91"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.83%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.56"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.75%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.17%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.356.98"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000268"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.29%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.401.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.40%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.665.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.77"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.90%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.80"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.77"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.525"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.36%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "530.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.12"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.46%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000204"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.54"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.85%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.18"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.870.19"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.53%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.31%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.47"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.71%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.17%  "

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +16.10%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.28%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "670.61"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +9.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.52%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "32.42"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.591"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.63%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.81"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.26%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.78%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.958"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.47%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.88%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "38.86"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +18.61%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0466"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.40%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.448"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +13.21%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.32"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.03%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.71"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.07%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.08%  "
